$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = 71.75571737916532
$ws.Range("D12").Value = 43.84675166390883
$ws.Range("D23").Value = 34.07741246683727
$ws.Range("D34").Value = 24.79350493092231
